$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates (refreshed crypto price/volume data).
# Numeric-looking text values need NumberFormat forced to Text first,
# otherwise Excel auto-converts them to the Number type; ClearFormats
# afterwards restores the original (unstyled) cell formatting.

$ws.Range("D2").Value = '63.021.25'
$ws.Range("E2").Value = '  +5.31%  '
$ws.Range("D3").Value = '3.369.32'
$ws.Range("E3").Value = '  +5.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +6.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.66'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.75%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.374.62'
$ws.Range("E8").Value = '  +5.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.42'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.118'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.34%  '
$ws.Range("E12").Value = '  +2.65%  '
$ws.Range("D13").Value = '3.945.44'
$ws.Range("E13").Value = '  +5.56%  '
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.07'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.50%  '
$ws.Range("E16").Value = '  +4.87%  '
$ws.Range("D17").Value = '63.064.99'
$ws.Range("E17").Value = '  +5.31%  '
$ws.Range("D18").Value = '3.330.98'
$ws.Range("E18").Value = '  +4.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.91'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.41'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.62'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +5.16%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.535'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.58'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.37'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +7.85%  '
$ws.Range("E27").Value = '  +6.57%  '
$ws.Range("D28").Value = '0.0₃0968'
$ws.Range("E28").Value = '  +10.53%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +6.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.01'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.58'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.31'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +10.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.27'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.70'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.55%  '
$ws.Range("E36").Value = '  +8.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.89'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.30%  '
$ws.Range("E38").Value = '  +12.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.94'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.98%  '
$ws.Range("D40").Value = '2.851.16'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("E41").Value = '  +5.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0326'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +10.45%  '
$ws.Range("E43").Value = '  +3.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.747'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.27'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("B46").Value = 'RenzoRestakedETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D46").Value = '3.411.25'
$ws.Range("E46").Value = '  +5.69%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.04'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +5.13%  '
$ws.Range("E48").Value = '  +7.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '298.95'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +13.45%  '
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.31'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.91%  '
